# Apply edits described by the diff:
#  1. Update the summary report title text (A1) with the new date range.
#  2. Update the start_date (col G) / end_date (col H) values for rows 5-24
#     (these hold shared strings "12/22/2019" / "1/4/2020" that become
#     "12/29/2019" / "1/11/2020").
#  3. Update the numeric "hours" values in column C for rows 5,6,7,10,11,
#     13-23 (the rest stay 0 and are left untouched).
#
# Dates are written through a Formula -> Copy -> PasteSpecial(values) round
# trip so the result lands back in the cell as plain text (shared string),
# matching the original cell type/format instead of being auto-parsed into
# a date serial number by a plain .Value assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

# 1. Title
$ws.Range("A1").Value = "Summary report for 12/29/2019 through 1/11/2020"

# 2. start_date / end_date columns (G = 7, H = 8) for rows 5-24
$newStart = "12/29/2019"
$newEnd = "1/11/2020"
for ($r = 5; $r -le 24; $r++) {
    Set-TextValue $ws.Cells.Item($r, 7) $newStart
    Set-TextValue $ws.Cells.Item($r, 8) $newEnd
}

# 3. hours column (C = 3)
$hours = @{
    5  = 4.9800000000000004
    6  = 45.02
    7  = 18.43
    10 = 24.48
    11 = 1.97
    13 = 28.27
    14 = 43.45
    15 = 66.7
    16 = 47.2
    17 = 8.8699999999999992
    18 = 18.93
    19 = 28.02
    20 = 40.42
    21 = 40.17
    22 = 42.73
    23 = 41.63
}

foreach ($r in $hours.Keys) {
    $ws.Cells.Item($r, 3).Value = $hours[$r]
}

# Clear the clipboard-marquee state left behind by Copy()
$excel.CutCopyMode = 0
